# Matriz de trazabilidad - actualizacion de estado
#
# Commit: "Recuperar contraseña, Autenticar y Registrar Usuario Tacna F&D
#          + Actualizacion Matriz de Trazabilidad"
#
# Marca como "Completado" los requerimientos "Registro Usuario" (fila 3),
# "Autenticar Usuario" (fila 4) y "Restablecer Contraseña de Usuario"
# (fila 6) -copiando el formato verde ya usado por las filas que estan
# "Completado"-, y actualiza la fecha de estado (columna I) del
# 2020-10-28 (44132) al 2020-10-31 (44135) para las filas que seguian
# con estado "Falta".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Columna H ("Estado"): Registro Usuario, Autenticar Usuario y
# Restablecer Contraseña de Usuario pasan de "Falta" (rojo) a
# "Completado" (verde). Se copia el formato de una fila ya completada
# (fila 16) para conservar el relleno/fuente verde usados en el resto
# de la matriz.
$completedRows = @(3, 4, 6)
foreach ($r in $completedRows) {
    $ws.Range("H16").Copy()
    $ws.Cells.Item($r, 8).PasteSpecial(-4122)
    $ws.Cells.Item($r, 8).Value = "Completado"
}
$excel.CutCopyMode = 0

# Columna I ("Fecha de estado"): se actualiza la fecha para las filas
# que mantienen el estado "Falta" (2020-10-28 -> 2020-10-31).
$dateRows = @(3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39)
foreach ($r in $dateRows) {
    $ws.Cells.Item($r, 9).Value = 44135
}

# Ultima seleccion activa del usuario antes de guardar.
$ws.Range("I30:I39").Select()
